# Autogen: updated excel table and exp output to match SR test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell value updates on the test table ---

# Row 5: fix correction of operator glyph used in the init call comment
$ws.Range("C5").Value = "Room(Run := TRUE, Temp:= 20.1); (* start temp *) "

# Row 8
$ws.Range("H8").Value = $false

# Row 9
$ws.Range("B9").Value = 2000
$ws.Range("H9").Value = $true

# Row 10
$ws.Range("B10").Value = 3000
$ws.Range("H10").Value = $false

# Row 12
$ws.Range("B12").Value = 0

# Row 13 unchanged

# Row 14
$ws.Range("B14").Value = 2000
$ws.Range("E14").Value = $true

# Row 15
$ws.Range("B15").Value = 3000
$ws.Range("J15").Value = "3000, 22.0"

# Trim trailing spaces in J12's text value
$ws.Range("J12").Value = "0, 22.1"

# --- View / selection updates ---
$ws.Range("C17").Select()

$wb.Windows.Item(1).Top = 4320
